$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 4 (Fabrica 3)
$ws.Range("B4").Value = "Hashlin Comelona"
$ws.Range("C4").Value = "India"
$ws.Range("D4").Value = 5000
$ws.Range("E4").Value = 3

# Update existing row 5 (Fabrica 4)
$ws.Range("B5").Value = "Serquen"
$ws.Range("C5").Value = "Cerru"
$ws.Range("D5").Value = 2000
$ws.Range("E5").Value = 30

# Add new row 6 (Fabrica 5), copying the style of A-column label cells (e.g. A5)
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A6").Value = "Fabrica 5:"
$ws.Range("B6").Value = "Remedial"
$ws.Range("C6").Value = "Cerru"
$ws.Range("D6").Value = 1300
$ws.Range("E6").Value = 30
